# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
# Rebuild the account-statement detail table (rows 16-50): instead of being
# grouped by mora period with every worker repeated per period, the data is
# now grouped by worker, listing their 5 mora periods together, and the
# "Salario Basico" column is normalized to a single value for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New roster of workers (N Doc Trabajador / Nombre Trabajador), in the order
# they now appear in the table.
$workers = @(
    @{ Doc = "1047421078"; Name = "JOSE ALFREDO ZUÑIGA GOMEZ" },
    @{ Doc = "1047398614"; Name = "MICHAEL JAVIER BUSTAMANTE BARRAGAN" },
    @{ Doc = "73214703";   Name = "JULIO CESAR JIMENEZ SILGADO" },
    @{ Doc = "1067094348"; Name = "MANUEL ANTONIO VASQUEZ BELLO" },
    @{ Doc = "1048292047"; Name = "YEISON MANUEL OLIVO CARDENAS" },
    @{ Doc = "80008129";   Name = "FERNEY PINTO MARIN" },
    @{ Doc = "1099549129"; Name = "JEFFERSON ALONSO MARTINEZ PARRA" }
)

# Mora periods, newest first, repeated for every worker.
$periods = @("1903", "1902", "1901", "1812", "1811")

$row = 16
foreach ($worker in $workers) {
    foreach ($period in $periods) {
        # Valor Mora: period 1903 carries 26041, every other period 31249.
        if ($period -eq "1903") {
            $valorMora = 26041
        } else {
            $valorMora = 31249
        }

        $ws.Cells.Item($row, 2).Value = "CC"
        $ws.Cells.Item($row, 3).Value = $worker.Doc
        $ws.Cells.Item($row, 4).Value = $worker.Name
        $ws.Cells.Item($row, 5).Value = $period
        $ws.Cells.Item($row, 6).Value = $valorMora
        $ws.Cells.Item($row, 7).Value = 781242

        $row = $row + 1
    }
}
